# Replace the 15 lattice-multiplication exercise cells (5 rows x 3 cols)
# in the single table with new problems, preserving the existing
# "NN x NN / digits / ---- / digit| | / digit| |" layout and run formatting.
$d = $word.ActiveDocument
$t = $d.Tables(1)
$nl = [char]11  # w:br maps to vertical-tab within Range.Text

$cell = $t.Cell(1, 1)
$cell.Range.Text = "81 x 48" + $nl + "  4    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"

$cell = $t.Cell(1, 2)
$cell.Range.Text = "61 x 73" + $nl + "  7    3" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"

$cell = $t.Cell(1, 3)
$cell.Range.Text = "54 x 22" + $nl + "  2    2" + $nl + "  ----" + $nl + "5|    |" + $nl + "4|    |"

$cell = $t.Cell(2, 1)
$cell.Range.Text = "33 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "3|    |" + $nl + "3|    |"

$cell = $t.Cell(2, 2)
$cell.Range.Text = "26 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"

$cell = $t.Cell(2, 3)
$cell.Range.Text = "30 x 25" + $nl + "  2    5" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"

$cell = $t.Cell(3, 1)
$cell.Range.Text = "85 x 87" + $nl + "  8    7" + $nl + "  ----" + $nl + "8|    |" + $nl + "5|    |"

$cell = $t.Cell(3, 2)
$cell.Range.Text = "50 x 23" + $nl + "  2    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "0|    |"

$cell = $t.Cell(3, 3)
$cell.Range.Text = "41 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "4|    |" + $nl + "1|    |"

$cell = $t.Cell(4, 1)
$cell.Range.Text = "59 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"

$cell = $t.Cell(4, 2)
$cell.Range.Text = "69 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"

$cell = $t.Cell(4, 3)
$cell.Range.Text = "77 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "7|    |" + $nl + "7|    |"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "67 x 54" + $nl + "  5    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"

$cell = $t.Cell(5, 2)
$cell.Range.Text = "32 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "3|    |" + $nl + "2|    |"

$cell = $t.Cell(5, 3)
$cell.Range.Text = "42 x 45" + $nl + "  4    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
